# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded in the status report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-02 13:16:10"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-02 13:15:58"
$wsZhCn.Range("K2").Value = "2016-09-02 13:16:32"

# de-de sheet: same two columns for the first file row. H2 mirrors the
# Overview sheet's "Latest HO Xliff Generate Date" value.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-02 13:16:10"
$wsDeDe.Range("K2").Value = "2016-09-02 13:16:39"
